$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 66 (2025-02-17, abs_activity): D and F change to 10 and 20
$ws.Range("D66").Value = 10
$ws.Range("F66").Value = 20

# Append new rows for 2025-02-18
$newRows = @(
    @("2025-02-18", "abs_activity", 9.791585277856878, 10, 0, 19.79158527785688),
    @("2025-02-18", "rel_activity", 0, 10, 0, 10),
    @("2025-02-18", "abs_sleep", 10, 10, 0, 20),
    @("2025-02-18", "rel_sleep", 10, 9.197440803635494, 0, 19.19744080363549)
)

$startRow = 70
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds a date-like string (e.g. "2025-02-18") which Excel would
    # otherwise auto-detect and convert to a real date value. Force it to be
    # stored as plain text, then restore the default (Normal) style so no
    # extra formatting is left behind on the cell.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
